# The document ends with a run of trailing empty paragraphs just before
# the section break (sectPr). Originally there were 4 of them; only the
# first should remain, so the last 3 (including the one holding two
# lone-space runs) are removed.
$d = $word.ActiveDocument

$count = $d.Paragraphs.Count

# Keep the first trailing empty paragraph (Paragraphs.Item($count - 2)),
# delete the remaining three up through the last paragraph of the body.
$firstToRemove = $d.Paragraphs.Item($count - 2)
$lastToRemove = $d.Paragraphs.Item($count)

$range = $d.Range($firstToRemove.Range.Start, $lastToRemove.Range.End)
$range.Delete()
